{"js": "// Rewrite the 10 \"Question N: {python-dict-repr}\" paragraphs into the\n// correct plain-text Q&A layout: each question paragraph followed by its\n// four lettered option paragraphs as separate paragraphs.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// New content, grouped per question: [questionText, optionA, optionB, optionC, optionD]\nconst groups = [\n  [\n    \"Question 1: A sample is composed of 2.78 g of iron and 1.19 g of oxygen. The empirical formula is ________.\",\n    \"A) FeO2\",\n    \"B) Fe2O5\",\n    \"C) Fe2O3\",\n    \"D) FeO\",\n  ],\n  [\n    \"Question 2: Which of the following compounds has its oxygen composition equal to 50.0%?\",\n    \"A) (NH4)2CO3\",\n    \"B) CH2O2\",\n    \"C) C2H3NO2\",\n    \"D) C10H20ONS\",\n  ],\n  [\n    \"Question 3: Which of the following compounds has its nitrogen composition equal to 29.2%?\",\n    \"A) C2H3NO2\",\n    \"B) C10H20ONS\",\n    \"C) (NH4)2CO3\",\n    \"D) N2H4\",\n  ],\n];\n\nconst originalParagraphs = paragraphs.items;\n\n// Reuse the first three existing paragraphs (old Question 1/2/3) as the new\n// question-header paragraphs, and insert the four option paragraphs right\n// after each one.\nfor (let i = 0; i < groups.length; i++) {\n  const [question, ...options] = groups[i];\n  const hostParagraph = originalParagraphs[i];\n  hostParagraph.insertText(question, Word.InsertLocation.replace);\n  let insertAfter = hostParagraph;\n  for (const option of options) {\n    insertAfter = insertAfter.insertParagraph(option, Word.InsertLocation.after);\n  }\n}\nawait context.sync();\n\n// Remove the leftover original question paragraphs (old Question 4..10).\nfor (let i = groups.length; i < originalParagraphs.length; i++) {\n  originalParagraphs[i].delete();\n}\nawait context.sync();\n", "ps1": "# Rewrite the 10 \"Question N: {python-dict-repr}\" paragraphs into the\n# correct plain-text Q&A layout: each question paragraph followed by its\n# four lettered option paragraphs as separate paragraphs.\n$d = $word.ActiveDocument\n\n# New content, grouped per question: question text, then options A-D.\n$groups = @(\n    @(\n        \"Question 1: A sample is composed of 2.78 g of iron and 1.19 g of oxygen. The empirical formula is ________.\",\n        \"A) FeO2\",\n        \"B) Fe2O5\",\n        \"C) Fe2O3\",\n        \"D) FeO\"\n    ),\n    @(\n        \"Question 2: Which of the following compounds has its oxygen composition equal to 50.0%?\",\n        \"A) (NH4)2CO3\",\n        \"B) CH2O2\",\n        \"C) C2H3NO2\",\n        \"D) C10H20ONS\"\n    ),\n    @(\n        \"Question 3: Which of the following compounds has its nitrogen composition equal to 29.2%?\",\n        \"A) C2H3NO2\",\n        \"B) C10H20ONS\",\n        \"C) (NH4)2CO3\",\n        \"D) N2H4\"\n    )\n)\n\n# Reuse the first three existing paragraphs (old Question 1/2/3) as the new\n# question-header paragraphs, and insert the four option paragraphs right\n# after each one.\n$paraIndex = 1\nforeach ($group in $groups) {\n    $questionText = $group[0]\n    $d.Paragraphs($paraIndex).Range.Text = $questionText\n    for ($i = 1; $i -lt $group.Count; $i++) {\n        $d.Paragraphs($paraIndex).Range.InsertParagraphAfter()\n        $paraIndex = $paraIndex + 1\n        $d.Paragraphs($paraIndex).Range.Text = $group[$i]\n    }\n    $paraIndex = $paraIndex + 1\n}\n\n# Remove the leftover original question paragraphs (old Question 4..10),\n# which now sit right after the newly written paragraphs. Deleting\n# repeatedly at $paraIndex removes them one by one since later paragraphs\n# shift down into that slot.\nwhile ($d.Paragraphs.Count -ge $paraIndex) {\n    $d.Paragraphs($paraIndex).Range.Delete()\n}\n"}
